# Actualización automática 2025-07-22 16:40:08
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L4").Value = 565.24
$ws1.Range("M4").Value = 7893.28
$ws1.Range("C25").Value = 1016.06
$ws1.Range("L56").Value = "3 de 54"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F4").Value = 10013.72
$ws2.Range("F25").Value = 11124.93
$ws2.Range("F56").Value = 55906.78

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D2").Value = 3089.66
$ws3.Range("E2").Value = 6880.68304517915
$ws3.Range("F2").Value = 0.3098850246174738

$ws3.Range("D15").Value = 3370.6
$ws3.Range("E15").Value = 10129.4
$ws3.Range("F15").Value = 0.2496740740740741

$ws3.Range("D16").Value = 46921.24
$ws3.Range("E16").Value = 4905.220000000001
$ws3.Range("F16").Value = 0.9053529799256982

$ws3.Range("D19").Value = 55906.78
$ws3.Range("E19").Value = 57799.67064517915
$ws3.Range("F19").Value = 0.4916764148628388
